$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the moment-of-inertia formula cells (L26:L29) from 0.5*1*1*1 to 1*1*1
$ws.Range("L26").Formula = "=1*1*1"
$ws.Range("L27:L29").Formula = "=1*1*1"

# Update the description text in B30 (was "I = 0.5 x M x R^2", now "I = M x R^2").
# Leading apostrophe preserves the cell's existing quote-prefix ("stored as text") style.
$ws.Range("B30").Value = "'I = M x R^2"

# Update the sheet view: scroll position and selection
$ws.Range("B30").Select()
